$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 361. This shifts the existing rows 361-461
# down to become rows 366-466.
$ws.Rows("361:365").Insert()

# Common (constant across the whole data table) column values.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$catId     = 100112020
$categoria = "Tomate"
$unidad    = "`$/bandeja 18 kilos"
$origen    = "Provincia de Limarí"
$kgUnidad  = 18
$clasif    = "Hortaliza"

# Data for the 5 newly inserted rows (row, fecha, variedad, calidad, volumen, precioMin, precioMax, precioProm, precioKg)
$newRows = @(
    @(361, 44588, "Larga vida", "Primera", 1600, 8500, 9000, 8750, 486),
    @(362, 44588, "Larga vida", "Segunda", 1300, 6500, 7000, 6750, 375),
    @(363, 44588, "Larga vida", "Tercera", 400,  4500, 5000, 4750, 264),
    @(364, 44588, "Semiduro",   "Primera", 1700, 5000, 5500, 5250, 292),
    @(365, 44588, "Semiduro",   "Segunda", 1100, 3000, 3500, 3250, 181)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2  = $mercadoId
    $ws.Cells.Item($r, 2).Value2  = $mercado
    $ws.Cells.Item($r, 3).Value2  = $region
    $ws.Cells.Item($r, 4).Value2  = $row[1]
    $ws.Cells.Item($r, 5).Value2  = $codreg
    $ws.Cells.Item($r, 6).Value2  = $catId
    $ws.Cells.Item($r, 7).Value2  = $categoria
    $ws.Cells.Item($r, 8).Value2  = $row[2]
    $ws.Cells.Item($r, 9).Value2  = $row[3]
    $ws.Cells.Item($r, 10).Value2 = $row[4]
    $ws.Cells.Item($r, 11).Value2 = $row[5]
    $ws.Cells.Item($r, 12).Value2 = $row[6]
    $ws.Cells.Item($r, 13).Value2 = $row[7]
    $ws.Cells.Item($r, 14).Value2 = $unidad
    $ws.Cells.Item($r, 15).Value2 = $origen
    $ws.Cells.Item($r, 16).Value2 = $row[8]
    $ws.Cells.Item($r, 17).Value2 = $kgUnidad
    $ws.Cells.Item($r, 18).Value2 = $clasif
}
